$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the weekly report entry in row 4 (new shared strings, in the order
#    they first appear: Task, Date, Description).
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "Creation of welcome, sign_up and login page"
$ws.Range("D4").Value = "14/05/2025"
$ws.Range("E4").Value = 'Creation of basic pages listed below the "Task" column but still to complete especially the styling aspect'

# Row 4 grows taller to fit the wrapped description text.
$ws.Rows("4:4").RowHeight = 42

# ---------------------------------------------------------------------------
# 2. Drop the now-empty helper cells in the spacer columns (F/G) that no
#    longer carry any formatting, while keeping F4 (it still carries the
#    separator fill/style).
# ---------------------------------------------------------------------------
$ws.Range("F2:G2").Clear()
$ws.Range("G4").Clear()
$ws.Range("F5:G12").Clear()

# ---------------------------------------------------------------------------
# 3. Apply "Wrap Text" across the whole report area. Each call only touches
#    cells that already shared one exact formatting (avoiding multi-area
#    references) so identical results get folded back into a single style,
#    mirroring Excel's own de-duplication of the style table when it
#    rewrites the file.
# ---------------------------------------------------------------------------
$ws.Range("C1").WrapText = $true
$ws.Range("D1:I1").WrapText = $true
$ws.Range("J1").WrapText = $true

$ws.Range("C2").WrapText = $true
$ws.Range("H2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("I2").WrapText = $true
$ws.Range("E2").WrapText = $true
$ws.Range("J2").WrapText = $true

$ws.Range("C3").WrapText = $true
$ws.Range("D3").WrapText = $true
$ws.Range("E3").WrapText = $true
$ws.Range("F3").WrapText = $true
$ws.Range("G3").WrapText = $true
$ws.Range("H3").WrapText = $true
$ws.Range("I3").WrapText = $true
$ws.Range("J3").WrapText = $true

$ws.Range("F4").WrapText = $true
$ws.Range("C4:C12").WrapText = $true
$ws.Range("D4:D12").WrapText = $true
$ws.Range("E4:E12").WrapText = $true
$ws.Range("H4:H12").WrapText = $true
$ws.Range("I4:I12").WrapText = $true
$ws.Range("J4:J12").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Leave the cursor where the author left it after typing the new row.
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()
